$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings are written as text, matching the
# original inlineStr cell type (Excel would otherwise coerce e.g. "16.20"
# into the number 16.2, silently dropping the trailing zero).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.818.08"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.40%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.342.59"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.29%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.32%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.59"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.26%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.667"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.54%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.02"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -6.40%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.09%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.589"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -7.44%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0993"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.45%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "32.35"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.29%  "

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.45%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.15"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -6.15%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.694.28"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.15%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.20"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -4.83%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.900"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.05%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.347.94"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.03%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.709.65"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.23%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000101"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.00%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.67"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.39%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "77.97"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.17%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.16"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.10%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.91"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +7.46%  "

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.11%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.49"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.73%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.35"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -8.15%  "

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.42%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "175.82"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.10%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.22"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.67%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.127"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.26%  "

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.09%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0736"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.49%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.08"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.26%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.36"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.51%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.73"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.02%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.37"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.14%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.37"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.43%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0271"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.02%  "

# Row 41
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "FTXToken"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.27"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +16.90%  "

# Row 42
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.01"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +18.73%  "

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.90%  "

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +5.42%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.85"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.33%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.196"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.86%  "

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.07%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.45"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.14%  "

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.43%  "

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.19%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "97.94"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.95%  "
